# Update cryptos list with latest price/volume data
# (mirrors the automated GitHub Actions refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.653.97'
$ws.Range("E2").Value = '  +1.56%  '
# Row 3
$ws.Range("D3").Value = '1.623.80'
$ws.Range("E3").Value = '  +2.14%  '
# Row 4
$ws.Range("E4").Value = '  -0.37%  '
# Row 5
$ws.Range("D5").Value = '''214.86'
$ws.Range("E5").Value = '  +1.23%  '
# Row 6
$ws.Range("D6").Value = '''0.506'
$ws.Range("E6").Value = '  +0.90%  '
# Row 7
$ws.Range("E7").Value = '  -0.32%  '
# Row 10
$ws.Range("E10").Value = '  +0.06%  '
# Row 11
$ws.Range("D11").Value = '''0.0857'
$ws.Range("E11").Value = '  +0.87%  '
# Row 12
$ws.Range("D12").Value = '1.849.64'
$ws.Range("E12").Value = '  +1.98%  '
# Row 13
$ws.Range("D13").Value = '1.631.26'
$ws.Range("E13").Value = '  +0.29%  '
# Row 14
$ws.Range("E14").Value = '  +0.36%  '
# Row 15
$ws.Range("D15").Value = '''0.514'
$ws.Range("E15").Value = '  -1.24%  '
# Row 16
$ws.Range("D16").Value = '''65.02'
$ws.Range("E16").Value = '  +1.01%  '
# Row 17
$ws.Range("D17").Value = '26.624.18'
$ws.Range("E17").Value = '  +1.37%  '
# Row 18
$ws.Range("D18").Value = '''232.03'
$ws.Range("E18").Value = '  +8.60%  '
# Row 20
$ws.Range("E20").Value = '  +2.88%  '
# Row 21
$ws.Range("E21").Value = '  -0.31%  '
# Row 22
$ws.Range("E22").Value = '  +2.40%  '
# Row 23
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").Value = '''2.24'
$ws.Range("E23").Value = '  +4.33%  '
# Row 24
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").Value = '''9.13'
$ws.Range("E24").Value = '  +1.25%  '
# Row 25
$ws.Range("D25").Value = '''145.52'
$ws.Range("E25").Value = '  +1.29%  '
# Row 26
$ws.Range("E26").Value = '  -0.13%  '
# Row 27
$ws.Range("D27").Value = '''7.06'
$ws.Range("E27").Value = '  +0.06%  '
# Row 28
$ws.Range("E28").Value = '  +2.35%  '
# Row 29
$ws.Range("D29").Value = '''15.67'
$ws.Range("E29").Value = '  +3.03%  '
# Row 30
$ws.Range("D30").Value = '''0.0500'
$ws.Range("E30").Value = '  +0.10%  '
# Row 31
$ws.Range("D31").Value = '''1.16'
$ws.Range("E31").Value = '  +0.69%  '
# Row 32
$ws.Range("E32").Value = '  +1.78%  '
# Row 33
$ws.Range("D33").Value = '1.446.99'
$ws.Range("E33").Value = '  +8.00%  '
# Row 34
$ws.Range("D34").Value = '''3.01'
$ws.Range("E34").Value = '  +2.45%  '
# Row 35
$ws.Range("E35").Value = '  -1.05%  '
# Row 36
$ws.Range("E36").Value = '  +0.74%  '
# Row 37
$ws.Range("D37").Value = '''0.560'
$ws.Range("E37").Value = '  -5.28%  '
# Row 38
$ws.Range("E38").Value = '  +0.67%  '
# Row 39
$ws.Range("D39").Value = '''0.840'
$ws.Range("E39").Value = '  +2.78%  '
# Row 40
$ws.Range("E40").Value = '  +2.06%  '
# Row 41
$ws.Range("E41").Value = '  -0.24%  '
# Row 42
$ws.Range("E42").Value = '  +3.08%  '
# Row 43
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '''0.943'
$ws.Range("E43").Value = '  -5.80%  '
# Row 44
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.761.71'
$ws.Range("E44").Value = '  +2.09%  '
# Row 45
$ws.Range("D45").Value = '''0.763'
$ws.Range("E45").Value = '  -0.34%  '
# Row 46
$ws.Range("D46").Value = '''62.23'
$ws.Range("E46").Value = '  +0.56%  '
# Row 47
$ws.Range("D47").Value = '''88.60'
$ws.Range("E47").Value = '  +3.37%  '
# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.51'
$ws.Range("E48").Value = '  +2.23%  '
# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  -1.37%  '
# Row 50
$ws.Range("E50").Value = '  +0.33%  '
# Row 51
$ws.Range("D51").Value = '''0.0969'
$ws.Range("E51").Value = '  -0.81%  '
